$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 75, shifting existing rows 75:113 down to 76:114
$ws.Rows("75:75").Insert()

# Populate the newly inserted row 75 with the new weekly record
$ws.Range("A75").Value = 7
$ws.Range("B75").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C75").Value = "Ñuble"
$ws.Range("D75").Value = 44438
$ws.Range("E75").Value = 16
$ws.Range("F75").Value = 100112003
$ws.Range("G75").Value = "Ajo"
$ws.Range("H75").Value = "Chino"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 40
$ws.Range("K75").Value = 14000
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = 14500
$ws.Range("N75").Value = "$/caja 10 kilos"
$ws.Range("O75").Value = "China"
$ws.Range("P75").Value = 1450
$ws.Range("Q75").Value = 10
$ws.Range("R75").Value = "Hortaliza"
